$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New row 14: a note below the existing "*Luu y" note row (row13).
# B14 reuses the red Times New Roman 13 font (no border/fill/alignment),
# C14:F14 reuse the plain black Times New Roman 13 font cells.
# ---------------------------------------------------------------------------
$ws.Rows(14).RowHeight = 17
$ws.Range("B14").Font.Name = "Times New Roman"
$ws.Range("B14").Font.Size = 13
$ws.Range("B14").Font.Color = 255
$ws.Range("B14").Value2 = "Số câu trong mỗi đề tối đa 120 câu"

$ws.Range("C14:F14").Font.Name = "Times New Roman"
$ws.Range("C14:F14").Font.Size = 13

# ---------------------------------------------------------------------------
# Row 10: A10 becomes a literal number 109 (no longer a shared text string).
# Some of the answer-key letters in this row are also edited.
# ---------------------------------------------------------------------------
$ws.Range("A10").Value2 = 109
$ws.Range("C10").Value2 = "B"
$ws.Range("E10").Value2 = "C"
$ws.Range("H10").Value2 = "D"
$ws.Range("I10").Value2 = "E"

# ---------------------------------------------------------------------------
# Row 2 (A2): was the literal number 101 (style s=2) -> becomes text "001"
# keeping the same visual style (font/fill/border/center) but with a
# Text ("@") number format, and the header label "Mã đề phải" row.
# ---------------------------------------------------------------------------
$ws.Range("B2").Copy() | Out-Null
$ws.Range("A2").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value2 = "001"

# ---------------------------------------------------------------------------
# Row 3 (A3): was "102" -> becomes text "002", style switches from the plain
# header style (s=1) to the same style but with Text ("@") number format.
# ---------------------------------------------------------------------------
$ws.Range("B1").Copy() | Out-Null
$ws.Range("A3").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value2 = "002"

# ---------------------------------------------------------------------------
# Row 4 (A4): was "103" -> becomes text "015", same style treatment as A3.
# ---------------------------------------------------------------------------
$ws.Range("B1").Copy() | Out-Null
$ws.Range("A4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value2 = "015"

# ---------------------------------------------------------------------------
# Rows 5-9 (A5:A9): their displayed text ("104".."108") is unchanged; only
# the shared-string table shifts because "102"/"103" are no longer used
# elsewhere. No direct edit is required here.
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# Update the view: drop the scrolled topLeftCell and move the active
# selection to H12.
# ---------------------------------------------------------------------------
$ws.Range("H12").Select() | Out-Null

Write-Output "edit complete"
